# edit.ps1 - applies the RQ1_Res.xlsx changes described by the commit
# "sootup getClasses func changed, tc24 25 cp bug, testing cbf"
#
# Summary of changes:
#  1. JB_CBF sheet (sheet11.xml): add a new header "cbf perform dead code
#     eliminator at the last step" in E1, and move the selection to E2.
#  2. JB_CP sheet (sheet7.xml):
#     - E33 gets a red font color applied (existing text, new style).
#     - Two new rows (35, 36) for tc24 / tc25 with a red-font comment in
#       column E ("not printed in parameters" / "Bug, should not
#       propagate inside array index, even in Soot").
#     - Five new rows (37-41) for tc26..tc30, each with "Y" in column C.
#     - Selection moves to C37.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) JB_CBF sheet - new column E header
# ---------------------------------------------------------------
$wsCbf = $wb.Worksheets.Item("JB_CBF")
$wsCbf.Activate() | Out-Null
$wsCbf.Range("E1").Value = "cbf perform dead code eliminator at the last step"
$wsCbf.Range("E2").Select() | Out-Null

# ---------------------------------------------------------------
# 2) JB_CP sheet - tc24 / tc25 comments + tc26-tc30 rows
# ---------------------------------------------------------------
$wsCp = $wb.Worksheets.Item("JB_CP")
$wsCp.Activate() | Out-Null

# Existing row 33 - comment in E33 becomes highlighted in red
$wsCp.Range("E33").Font.Color = 255

# New row 35 - tc24
$wsCp.Range("A35").Value = "tc24"
$wsCp.Range("E35").Value = "not printed in parameters"
$wsCp.Range("E35").Font.Color = 255

# New row 36 - tc25
$wsCp.Range("A36").Value = "tc25"
$wsCp.Range("E36").Value = "Bug, should not propagate inside array index, even in Soot"
$wsCp.Range("E36").Font.Color = 255

# New rows 37-41 - tc26..tc30, all marked "Y" in column C
$wsCp.Range("A37").Value = "tc26"
$wsCp.Range("C37").Value = "Y"

$wsCp.Range("A38").Value = "tc27"
$wsCp.Range("C38").Value = "Y"

$wsCp.Range("A39").Value = "tc28"
$wsCp.Range("C39").Value = "Y"

$wsCp.Range("A40").Value = "tc29"
$wsCp.Range("C40").Value = "Y"

$wsCp.Range("A41").Value = "tc30"
$wsCp.Range("C41").Value = "Y"

# Final selection as left by the author
$wsCp.Range("C37").Select() | Out-Null
